# The commit swaps the presentation's applied design from the custom
# "Integral" theme back to the stock "Office Theme" palette (the
# "Integral" colors end up parked on the theme part that the Notes
# Master points at instead). The slide master / color scheme is the
# part of that swap that's reachable through the PowerPoint object
# model, so recolor the deck's theme color scheme (12 slots, in the
# standard dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink order) from the
# "Integral" palette to the stock "Office Theme" palette.

$p = $ppt.ActivePresentation

$design = $p.Designs.Item(1)

# Best-effort: real PowerPoint exposes Design.Name as read-only too, but
# try it in case a given host allows the relabel.
try { $design.Name = "Office Theme" } catch { }

$master = $design.SlideMaster
$theme  = $master.Theme
$colors = $theme.ThemeColorScheme

# Target "Office Theme" color scheme, expressed as COM RGB() long values
# (0xBBGGRR) so they round-trip through ColorFormat.RGB.
#   1  dk1      #000000
#   2  lt1      #FFFFFF
#   3  dk2      #44546A
#   4  lt2      #E7E6E6
#   5  accent1  #5B9BD5
#   6  accent2  #ED7D31
#   7  accent3  #A5A5A5
#   8  accent4  #FFC000
#   9  accent5  #4472C4
#   10 accent6  #70AD47
#   11 hlink    #0563C1
#   12 folHlink #954F72
$officeThemeRgb = @(
    0,
    16777215,
    6968388,
    15132391,
    13998939,
    3243501,
    10855845,
    49407,
    12874308,
    4697456,
    12673797,
    7491477
)

for ($i = 1; $i -le $colors.Count; $i++) {
    $colors.Item($i).RGB = $officeThemeRgb[$i - 1]
}
